$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "human"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "mouse"

$ws2.Range("A1").Value = "sample"
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Color = 0x000000
$ws2.Range("A1").Borders.Item(9).LineStyle = 1
$ws2.Range("A1").Borders.Item(9).ColorIndex = 64
$ws2.Range("A1").WrapText = $true
